# Applies per-row Price (D) / Volume(1h) (E) updates, and the
# Mantle/OKB row swap (B/C/D/E for rows 45-46), to match the
# "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = '61.007.97'
$ws.Cells.Item(2,5).Value = '  -2.13%  '

$ws.Cells.Item(3,4).Value = '3.374.11'
$ws.Cells.Item(3,5).Value = '  +0.15%  '

$ws.Cells.Item(4,5).Value = '  +0.00%  '

$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = '572.17'
$ws.Cells.Item(5,5).Value = '  +1.18%  '

$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = '135.65'
$ws.Cells.Item(6,5).Value = '  +9.33%  '

$ws.Cells.Item(7,5).Value = '  -0.03%  '

$ws.Cells.Item(8,4).Value = '3.373.48'
$ws.Cells.Item(8,5).Value = '  +0.17%  '

$ws.Cells.Item(9,5).Value = '  +1.45%  '

$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = '7.60'
$ws.Cells.Item(10,5).Value = '  +5.48%  '

$ws.Cells.Item(11,5).Value = '  +3.52%  '

$ws.Cells.Item(12,5).Value = '  +4.27%  '

$ws.Cells.Item(13,4).Value = '3.951.06'
$ws.Cells.Item(13,5).Value = '  +0.31%  '

$ws.Cells.Item(14,5).Value = '  +1.71%  '

$ws.Cells.Item(15,4).NumberFormat = "@"
$ws.Cells.Item(15,4).Value = '0.0000173'
$ws.Cells.Item(15,5).Value = '  +2.38%  '

$ws.Cells.Item(16,4).Value = '3.380.00'
$ws.Cells.Item(16,5).Value = '  +0.32%  '

$ws.Cells.Item(17,4).NumberFormat = "@"
$ws.Cells.Item(17,4).Value = '25.10'
$ws.Cells.Item(17,5).Value = '  +3.51%  '

$ws.Cells.Item(18,4).Value = '61.179.13'
$ws.Cells.Item(18,5).Value = '  -2.02%  '

$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(19,4).Value = '14.03'
$ws.Cells.Item(19,5).Value = '  +8.05%  '

$ws.Cells.Item(20,5).Value = '  +3.86%  '

$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,4).Value = '9.39'
$ws.Cells.Item(21,5).Value = '  +1.73%  '

$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value = '374.05'
$ws.Cells.Item(22,5).Value = '  +1.70%  '

$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value = '0.568'
$ws.Cells.Item(23,5).Value = '  +3.12%  '

$ws.Cells.Item(24,4).Value = '3.509.84'
$ws.Cells.Item(24,5).Value = '  +0.17%  '

$ws.Cells.Item(25,5).Value = '  +0.06%  '

$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(26,4).Value = '70.60'
$ws.Cells.Item(26,5).Value = '  -0.29%  '

$ws.Cells.Item(27,5).Value = '  +12.54%  '

$ws.Cells.Item(28,5).Value = '  +22.29%  '

$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(29,4).Value = '7.68'
$ws.Cells.Item(29,5).Value = '  +13.32%  '

$ws.Cells.Item(30,5).Value = '  +0.02%  '

$ws.Cells.Item(31,4).NumberFormat = "@"
$ws.Cells.Item(31,4).Value = '8.10'
$ws.Cells.Item(31,5).Value = '  +5.55%  '

$ws.Cells.Item(32,5).Value = '  +2.63%  '

$ws.Cells.Item(33,4).NumberFormat = "@"
$ws.Cells.Item(33,4).Value = '0.154'
$ws.Cells.Item(33,5).Value = '  +4.48%  '

$ws.Cells.Item(34,5).Value = '  -0.06%  '

$ws.Cells.Item(35,4).Value = '3.405.92'
$ws.Cells.Item(35,5).Value = '  +0.25%  '

$ws.Cells.Item(36,4).NumberFormat = "@"
$ws.Cells.Item(36,4).Value = '23.38'
$ws.Cells.Item(36,5).Value = '  +4.00%  '

$ws.Cells.Item(37,4).NumberFormat = "@"
$ws.Cells.Item(37,4).Value = '5.55'
$ws.Cells.Item(37,5).Value = '  +9.94%  '

$ws.Cells.Item(38,5).Value = '  +7.07%  '

$ws.Cells.Item(39,5).Value = '  +5.18%  '

$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value = '162.85'
$ws.Cells.Item(40,5).Value = '  -0.86%  '

$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,4).Value = '0.0789'
$ws.Cells.Item(41,5).Value = '  +6.19%  '

$ws.Cells.Item(42,5).Value = '  +0.04%  '

$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value = '4.42'
$ws.Cells.Item(43,5).Value = '  +5.31%  '

$ws.Cells.Item(44,4).NumberFormat = "@"
$ws.Cells.Item(44,4).Value = '1.20'
$ws.Cells.Item(44,5).Value = '  +14.37%  '

$ws.Cells.Item(45,2).Value = 'OKB'
$ws.Cells.Item(45,3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value = '41.38'
$ws.Cells.Item(45,5).Value = '  +0.90%  '

$ws.Cells.Item(46,2).Value = 'Mantle'
$ws.Cells.Item(46,3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,4).Value = '0.760'
$ws.Cells.Item(46,5).Value = '  -0.30%  '

$ws.Cells.Item(47,5).Value = '  +5.52%  '

$ws.Cells.Item(48,4).NumberFormat = "@"
$ws.Cells.Item(48,4).Value = '23.14'
$ws.Cells.Item(48,5).Value = '  +3.67%  '

$ws.Cells.Item(49,4).NumberFormat = "@"
$ws.Cells.Item(49,4).Value = '6.99'
$ws.Cells.Item(49,5).Value = '  +6.38%  '

$ws.Cells.Item(50,4).NumberFormat = "@"
$ws.Cells.Item(50,4).Value = '22.97'
$ws.Cells.Item(50,5).Value = '  +15.38%  '

$ws.Cells.Item(51,4).NumberFormat = "@"
$ws.Cells.Item(51,4).Value = '0.895'
$ws.Cells.Item(51,5).Value = '  +7.14%  '

Write-Host "Applied all changes"
